$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, shifting existing rows (111-201) down to (112-202)
$ws.Rows("111:111").Insert()

# Populate the newly inserted row 111 with the new record's data
$ws.Cells.Item(111, 1).Value = 3
$ws.Cells.Item(111, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44447
$ws.Cells.Item(111, 5).Value = 5
$ws.Cells.Item(111, 6).Value = 100112032
$ws.Cells.Item(111, 7).Value = "Zapallo italiano"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 85
$ws.Cells.Item(111, 11).Value = 13000
$ws.Cells.Item(111, 12).Value = 13500
$ws.Cells.Item(111, 13).Value = 13265
$ws.Cells.Item(111, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 190
$ws.Cells.Item(111, 17).Value = 70
$ws.Cells.Item(111, 18).Value = "Hortaliza"
